$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 25
$ws.Cells.Item(3, 3).Value = 45
$ws.Cells.Item(4, 3).Value = 84
$ws.Cells.Item(5, 3).Value = 71
$ws.Cells.Item(6, 3).Value = 108
$ws.Cells.Item(7, 3).Value = 73
$ws.Cells.Item(8, 3).Value = 46
$ws.Cells.Item(9, 3).Value = 51
$ws.Cells.Item(10, 3).Value = 28
$ws.Cells.Item(11, 3).Value = 63
$ws.Cells.Item(12, 3).Value = 37
$ws.Cells.Item(13, 3).Value = 50
$ws.Cells.Item(14, 3).Value = 95
$ws.Cells.Item(15, 3).Value = 59
$ws.Cells.Item(16, 3).Value = 74
$ws.Cells.Item(17, 3).Value = 34
$ws.Cells.Item(18, 3).Value = 44
$ws.Cells.Item(19, 3).Value = 47
$ws.Cells.Item(20, 3).Value = 32
$ws.Cells.Item(21, 3).Value = 64
$ws.Cells.Item(22, 3).Value = 34
$ws.Cells.Item(23, 3).Value = 25
$ws.Cells.Item(24, 3).Value = 67
$ws.Cells.Item(26, 3).Value = 39
$ws.Cells.Item(27, 3).Value = 57
$ws.Cells.Item(28, 3).Value = 56
$ws.Cells.Item(29, 3).Value = 69
$ws.Cells.Item(30, 3).Value = 56
$ws.Cells.Item(31, 3).Value = 61
$ws.Cells.Item(32, 3).Value = 44
$ws.Cells.Item(33, 3).Value = 49
$ws.Cells.Item(34, 3).Value = 60
$ws.Cells.Item(35, 3).Value = 56
$ws.Cells.Item(36, 3).Value = 39
$ws.Cells.Item(37, 3).Value = 58
$ws.Cells.Item(38, 3).Value = 64
$ws.Cells.Item(39, 3).Value = 62
$ws.Cells.Item(40, 3).Value = 78
$ws.Cells.Item(41, 3).Value = 60
$ws.Cells.Item(42, 3).Value = 75
$ws.Cells.Item(43, 3).Value = 62
$ws.Cells.Item(44, 3).Value = 67
$ws.Cells.Item(45, 3).Value = 38
$ws.Cells.Item(46, 3).Value = 67
$ws.Cells.Item(47, 3).Value = 53
$ws.Cells.Item(48, 3).Value = 51
$ws.Cells.Item(49, 3).Value = 60
$ws.Cells.Item(50, 3).Value = 41
$ws.Cells.Item(52, 3).Value = 49
$ws.Cells.Item(53, 3).Value = 83
$ws.Cells.Item(54, 3).Value = 53
$ws.Cells.Item(55, 3).Value = 65
$ws.Cells.Item(56, 3).Value = 88
$ws.Cells.Item(57, 3).Value = 69
$ws.Cells.Item(58, 3).Value = 71
$ws.Cells.Item(59, 3).Value = 62
$ws.Cells.Item(60, 3).Value = 54
$ws.Cells.Item(61, 3).Value = 37
$ws.Cells.Item(62, 3).Value = 62
$ws.Cells.Item(63, 3).Value = 105
$ws.Cells.Item(65, 3).Value = 32
$ws.Cells.Item(66, 3).Value = 63
$ws.Cells.Item(67, 3).Value = 80
$ws.Cells.Item(68, 3).Value = 55
$ws.Cells.Item(69, 3).Value = 96
$ws.Cells.Item(70, 3).Value = 51
$ws.Cells.Item(71, 3).Value = 71
$ws.Cells.Item(72, 3).Value = 41
$ws.Cells.Item(77, 3).Value = 131
$ws.Cells.Item(78, 3).Value = 78
$ws.Cells.Item(79, 3).Value = 109
$ws.Cells.Item(80, 3).Value = 109
$ws.Cells.Item(81, 3).Value = 73
$ws.Cells.Item(82, 3).Value = 27
$ws.Cells.Item(83, 3).Value = 70
$ws.Cells.Item(84, 3).Value = 186
$ws.Cells.Item(86, 3).Value = 13
$ws.Cells.Item(92, 3).Value = 241
$ws.Cells.Item(93, 3).Value = 5410
$ws.Cells.Item(2, 4).Value = 22
$ws.Cells.Item(3, 4).Value = 39
$ws.Cells.Item(4, 4).Value = 67
$ws.Cells.Item(5, 4).Value = 65
$ws.Cells.Item(6, 4).Value = 88
$ws.Cells.Item(7, 4).Value = 56
$ws.Cells.Item(8, 4).Value = 30
$ws.Cells.Item(9, 4).Value = 43
$ws.Cells.Item(11, 4).Value = 53
$ws.Cells.Item(12, 4).Value = 32
$ws.Cells.Item(13, 4).Value = 36
$ws.Cells.Item(14, 4).Value = 84
$ws.Cells.Item(15, 4).Value = 47
$ws.Cells.Item(16, 4).Value = 59
$ws.Cells.Item(17, 4).Value = 28
$ws.Cells.Item(18, 4).Value = 38
$ws.Cells.Item(20, 4).Value = 25
$ws.Cells.Item(21, 4).Value = 51
$ws.Cells.Item(22, 4).Value = 25
$ws.Cells.Item(23, 4).Value = 20
$ws.Cells.Item(24, 4).Value = 62
$ws.Cells.Item(25, 4).Value = 30
$ws.Cells.Item(26, 4).Value = 33
$ws.Cells.Item(27, 4).Value = 43
$ws.Cells.Item(28, 4).Value = 43
$ws.Cells.Item(29, 4).Value = 53
$ws.Cells.Item(30, 4).Value = 44
$ws.Cells.Item(31, 4).Value = 47
$ws.Cells.Item(32, 4).Value = 35
$ws.Cells.Item(33, 4).Value = 35
$ws.Cells.Item(34, 4).Value = 41
$ws.Cells.Item(35, 4).Value = 48
$ws.Cells.Item(36, 4).Value = 25
$ws.Cells.Item(37, 4).Value = 41
$ws.Cells.Item(38, 4).Value = 53
$ws.Cells.Item(39, 4).Value = 56
$ws.Cells.Item(40, 4).Value = 68
$ws.Cells.Item(41, 4).Value = 50
$ws.Cells.Item(42, 4).Value = 62
$ws.Cells.Item(43, 4).Value = 52
$ws.Cells.Item(44, 4).Value = 56
$ws.Cells.Item(45, 4).Value = 29
$ws.Cells.Item(46, 4).Value = 51
$ws.Cells.Item(47, 4).Value = 47
$ws.Cells.Item(48, 4).Value = 34
$ws.Cells.Item(49, 4).Value = 52
$ws.Cells.Item(50, 4).Value = 33
$ws.Cells.Item(51, 4).Value = 30
$ws.Cells.Item(52, 4).Value = 40
$ws.Cells.Item(53, 4).Value = 74
$ws.Cells.Item(54, 4).Value = 32
$ws.Cells.Item(55, 4).Value = 45
$ws.Cells.Item(56, 4).Value = 69
$ws.Cells.Item(57, 4).Value = 45
$ws.Cells.Item(58, 4).Value = 58
$ws.Cells.Item(59, 4).Value = 49
$ws.Cells.Item(60, 4).Value = 37
$ws.Cells.Item(61, 4).Value = 31
$ws.Cells.Item(62, 4).Value = 52
$ws.Cells.Item(63, 4).Value = 96
$ws.Cells.Item(65, 4).Value = 27
$ws.Cells.Item(66, 4).Value = 48
$ws.Cells.Item(67, 4).Value = 65
$ws.Cells.Item(68, 4).Value = 43
$ws.Cells.Item(69, 4).Value = 79
$ws.Cells.Item(70, 4).Value = 42
$ws.Cells.Item(71, 4).Value = 64
$ws.Cells.Item(72, 4).Value = 24
$ws.Cells.Item(73, 4).Value = 39
$ws.Cells.Item(77, 4).Value = 130
$ws.Cells.Item(78, 4).Value = 61
$ws.Cells.Item(79, 4).Value = 109
$ws.Cells.Item(80, 4).Value = 91
$ws.Cells.Item(81, 4).Value = 53
$ws.Cells.Item(83, 4).Value = 62
$ws.Cells.Item(84, 4).Value = 124
$ws.Cells.Item(86, 4).Value = 11
$ws.Cells.Item(92, 4).Value = 172
$ws.Cells.Item(93, 4).Value = 4344
